$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 746.7143
$ws.Cells.Item(32, 9).Value = 607.2
$ws.Cells.Item(32, 10).Value = 790.3125
$ws.Cells.Item(32, 11).Value = 607.2
$ws.Cells.Item(32, 12).Value = 790.3125
$ws.Cells.Item(32, 13).Value = -281.2
$ws.Cells.Item(32, 14).Value = -1442.3125
$ws.Cells.Item(40, 8).Value = 1317.2667
$ws.Cells.Item(40, 9).Value = 999
$ws.Cells.Item(40, 10).Value = 1476.4
$ws.Cells.Item(40, 11).Value = 999
$ws.Cells.Item(40, 12).Value = 1476.4
$ws.Cells.Item(40, 13).Value = -824
$ws.Cells.Item(40, 14).Value = -1826.4
$ws.Cells.Item(51, 8).Value = 8011.2
$ws.Cells.Item(51, 10).Value = 8050
$ws.Cells.Item(51, 12).Value = 8050
$ws.Cells.Item(51, 14).Value = -9018
$ws.Cells.Item(93, 8).Value = 21173.043
$ws.Cells.Item(93, 10).Value = 21173.043
$ws.Cells.Item(93, 12).Value = 21173.043
$ws.Cells.Item(93, 14).Value = -26165.043
$ws.Cells.Item(98, 8).Value = 6324.36
$ws.Cells.Item(98, 9).Value = 4650.5
$ws.Cells.Item(98, 10).Value = 7440.2666
$ws.Cells.Item(98, 11).Value = 4650.5
$ws.Cells.Item(98, 12).Value = 7440.2666
$ws.Cells.Item(98, 13).Value = -3152.5
$ws.Cells.Item(98, 14).Value = -10436.2666
$ws.Cells.Item(103, 8).Value = 6109.625
$ws.Cells.Item(103, 9).Value = 797.63635
$ws.Cells.Item(103, 10).Value = 17796
$ws.Cells.Item(103, 11).Value = 2392.90905
$ws.Cells.Item(103, 12).Value = 53388
$ws.Cells.Item(103, 13).Value = -1806.90905
$ws.Cells.Item(103, 14).Value = -54560
$ws.Cells.Item(113, 8).Value = 12383.286
$ws.Cells.Item(113, 10).Value = 20000
$ws.Cells.Item(113, 12).Value = 20000
$ws.Cells.Item(113, 14).Value = -26508
$ws.Cells.Item(116, 8).Value = 459166.72
$ws.Cells.Item(116, 9).Value = 835139
$ws.Cells.Item(116, 10).Value = 8000
$ws.Cells.Item(116, 11).Value = 835139
$ws.Cells.Item(116, 12).Value = 8000
$ws.Cells.Item(116, 13).Value = -831697
$ws.Cells.Item(116, 14).Value = -14884
$ws.Cells.Item(122, 8).Value = 6324.36
$ws.Cells.Item(122, 9).Value = 4650.5
$ws.Cells.Item(122, 10).Value = 7440.2666
$ws.Cells.Item(122, 11).Value = 13951.5
$ws.Cells.Item(122, 12).Value = 22320.7998
$ws.Cells.Item(122, 13).Value = -11501.5
$ws.Cells.Item(122, 14).Value = -27220.7998
$ws.Cells.Item(132, 8).Value = 27406262
$ws.Cells.Item(132, 9).Value = 34618750
$ws.Cells.Item(132, 11).Value = 103856250
$ws.Cells.Item(132, 13).Value = -103853720
$ws.Cells.Item(137, 8).Value = 868926.9399999999
$ws.Cells.Item(137, 9).Value = 1645061.9
$ws.Cells.Item(137, 11).Value = 4935185.699999999
$ws.Cells.Item(137, 13).Value = -4932635.699999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1280.7646
$ws.Cells.Item(2, 9).Value = 1291.4828
$ws.Cells.Item(2, 10).Value = 1218.6
$ws.Cells.Item(2, 11).Value = 1291.4828
$ws.Cells.Item(2, 12).Value = 1218.6
$ws.Cells.Item(2, 13).Value = -1178.4828
$ws.Cells.Item(2, 14).Value = -1444.6
$ws.Cells.Item(32, 8).Value = 5961.551
$ws.Cells.Item(32, 9).Value = 7232.5483
$ws.Cells.Item(32, 10).Value = 3772.611
$ws.Cells.Item(32, 11).Value = 7232.5483
$ws.Cells.Item(32, 12).Value = 3772.611
$ws.Cells.Item(32, 13).Value = -6945.5483
$ws.Cells.Item(32, 14).Value = -4346.611
$ws.Cells.Item(97, 8).Value = 513.94116
$ws.Cells.Item(97, 9).Value = 409.7857
$ws.Cells.Item(97, 11).Value = 409.7857
$ws.Cells.Item(97, 13).Value = 86.21429999999998
$ws.Cells.Item(103, 8).Value = 34083.332
$ws.Cells.Item(103, 10).Value = 34083.332
$ws.Cells.Item(103, 12).Value = 34083.332
$ws.Cells.Item(103, 14).Value = -36427.332
$ws.Cells.Item(116, 8).Value = 1280.7646
$ws.Cells.Item(116, 9).Value = 1291.4828
$ws.Cells.Item(116, 10).Value = 1218.6
$ws.Cells.Item(116, 11).Value = 1291.4828
$ws.Cells.Item(116, 12).Value = 1218.6
$ws.Cells.Item(116, 13).Value = 1002.5172
$ws.Cells.Item(116, 14).Value = -5806.6
$ws.Cells.Item(128, 8).Value = 41880
$ws.Cells.Item(128, 10).Value = 41880
$ws.Cells.Item(128, 12).Value = 41880
$ws.Cells.Item(128, 14).Value = -51840
$ws.Cells.Item(137, 8).Value = 48307.6
$ws.Cells.Item(137, 10).Value = 48307.6
$ws.Cells.Item(137, 12).Value = 48307.6
$ws.Cells.Item(137, 14).Value = -58507.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1280.7646
$ws.Cells.Item(3, 9).Value = 1291.4828
$ws.Cells.Item(3, 10).Value = 1218.6
$ws.Cells.Item(3, 11).Value = 1291.4828
$ws.Cells.Item(3, 12).Value = 1218.6
$ws.Cells.Item(3, 13).Value = -1177.4828
$ws.Cells.Item(3, 14).Value = -1446.6
$ws.Cells.Item(94, 8).Value = 557.4516
$ws.Cells.Item(94, 9).Value = 466.95834
$ws.Cells.Item(94, 10).Value = 867.7143
$ws.Cells.Item(94, 11).Value = 466.95834
$ws.Cells.Item(94, 12).Value = 867.7143
$ws.Cells.Item(94, 13).Value = -15.95834000000002
$ws.Cells.Item(94, 14).Value = -1769.7143
$ws.Cells.Item(95, 8).Value = 32090.908
$ws.Cells.Item(95, 10).Value = 32090.908
$ws.Cells.Item(95, 12).Value = 32090.908
$ws.Cells.Item(95, 14).Value = -37582.908
$ws.Cells.Item(103, 8).Value = 34090.91
$ws.Cells.Item(103, 10).Value = 34090.91
$ws.Cells.Item(103, 12).Value = 34090.91
$ws.Cells.Item(103, 14).Value = -36434.91
$ws.Cells.Item(105, 8).Value = 1733.1409
$ws.Cells.Item(105, 9).Value = 1677.0952
$ws.Cells.Item(105, 10).Value = 2174.5
$ws.Cells.Item(105, 11).Value = 1677.0952
$ws.Cells.Item(105, 12).Value = 2174.5
$ws.Cells.Item(105, 13).Value = 69.90480000000002
$ws.Cells.Item(105, 14).Value = -5668.5
$ws.Cells.Item(129, 8).Value = 43684.223
$ws.Cells.Item(129, 10).Value = 43684.223
$ws.Cells.Item(129, 12).Value = 43684.223
$ws.Cells.Item(129, 14).Value = -53684.223
$ws.Cells.Item(134, 8).Value = 3389.4583
$ws.Cells.Item(134, 9).Value = 1356.9166
$ws.Cells.Item(134, 10).Value = 5422
$ws.Cells.Item(134, 11).Value = 4070.7498
$ws.Cells.Item(134, 12).Value = 16266
$ws.Cells.Item(134, 13).Value = -1535.7498
$ws.Cells.Item(134, 14).Value = -21336
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 252592.64
$ws.Cells.Item(31, 9).Value = 615256.25
$ws.Cells.Item(31, 10).Value = 3261.4062
$ws.Cells.Item(31, 11).Value = 615256.25
$ws.Cells.Item(31, 12).Value = 3261.4062
$ws.Cells.Item(31, 13).Value = -614961.25
$ws.Cells.Item(31, 14).Value = -3851.4062
$ws.Cells.Item(34, 8).Value = 252592.64
$ws.Cells.Item(34, 9).Value = 615256.25
$ws.Cells.Item(34, 10).Value = 3261.4062
$ws.Cells.Item(34, 11).Value = 615256.25
$ws.Cells.Item(34, 12).Value = 3261.4062
$ws.Cells.Item(34, 13).Value = -615054.25
$ws.Cells.Item(34, 14).Value = -3665.4062
$ws.Cells.Item(137, 8).Value = 43380
$ws.Cells.Item(137, 10).Value = 43380
$ws.Cells.Item(137, 12).Value = 43380
$ws.Cells.Item(137, 14).Value = -53580
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 152.8
$ws.Cells.Item(33, 9).Value = 135.4
$ws.Cells.Item(33, 11).Value = 812.4000000000001
$ws.Cells.Item(33, 13).Value = -529.4000000000001
$ws.Cells.Item(68, 8).Value = 2384.739
$ws.Cells.Item(68, 9).Value = 659
$ws.Cells.Item(68, 10).Value = 4380.125
$ws.Cells.Item(68, 11).Value = 1977
$ws.Cells.Item(68, 12).Value = 13140.375
$ws.Cells.Item(68, 13).Value = -1166
$ws.Cells.Item(68, 14).Value = -14762.375
$ws.Cells.Item(71, 8).Value = 2384.739
$ws.Cells.Item(71, 9).Value = 659
$ws.Cells.Item(71, 10).Value = 4380.125
$ws.Cells.Item(71, 11).Value = 5931
$ws.Cells.Item(71, 12).Value = 39421.125
$ws.Cells.Item(71, 13).Value = -1875
$ws.Cells.Item(71, 14).Value = -47533.125
$ws.Cells.Item(92, 8).Value = 1506.2
$ws.Cells.Item(92, 9).Value = 1500
$ws.Cells.Item(92, 10).Value = 1506.8889
$ws.Cells.Item(92, 11).Value = 4500
$ws.Cells.Item(92, 12).Value = 4520.6667
$ws.Cells.Item(92, 13).Value = -3252
$ws.Cells.Item(92, 14).Value = -7016.6667
$ws.Cells.Item(129, 8).Value = 1386.64
$ws.Cells.Item(129, 9).Value = 1103.8462
$ws.Cells.Item(129, 10).Value = 1693
$ws.Cells.Item(129, 11).Value = 3311.5386
$ws.Cells.Item(129, 12).Value = 5079
$ws.Cells.Item(129, 13).Value = 1688.4614
$ws.Cells.Item(129, 14).Value = -15079
$ws.Cells.Item(131, 8).Value = 779.95
$ws.Cells.Item(131, 10).Value = 805.2283
$ws.Cells.Item(131, 12).Value = 2415.6849
$ws.Cells.Item(131, 14).Value = -12495.6849
$ws.Cells.Item(132, 8).Value = 4055.6428
$ws.Cells.Item(132, 9).Value = 616.4
$ws.Cells.Item(132, 11).Value = 5547.599999999999
$ws.Cells.Item(132, 13).Value = -3017.599999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 50003840
$ws.Cells.Item(80, 10).Value = 5733.3335
$ws.Cells.Item(80, 12).Value = 5733.3335
$ws.Cells.Item(80, 14).Value = -7729.3335
$ws.Cells.Item(83, 8).Value = 50003840
$ws.Cells.Item(83, 10).Value = 5733.3335
$ws.Cells.Item(83, 12).Value = 28666.6675
$ws.Cells.Item(83, 14).Value = -38650.6675
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 93954.63
$ws.Cells.Item(22, 9).Value = 251375.25
$ws.Cells.Item(22, 10).Value = 4000
$ws.Cells.Item(22, 11).Value = 251375.25
$ws.Cells.Item(22, 12).Value = 4000
$ws.Cells.Item(22, 13).Value = -251080.25
$ws.Cells.Item(22, 14).Value = -4590
$ws.Cells.Item(27, 8).Value = 93954.63
$ws.Cells.Item(27, 9).Value = 251375.25
$ws.Cells.Item(27, 10).Value = 4000
$ws.Cells.Item(27, 11).Value = 251375.25
$ws.Cells.Item(27, 12).Value = 4000
$ws.Cells.Item(27, 13).Value = -251268.25
$ws.Cells.Item(27, 14).Value = -4214
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 70607.336
$ws.Cells.Item(46, 10).Value = 70607.336
$ws.Cells.Item(46, 12).Value = 70607.336
$ws.Cells.Item(46, 14).Value = -71069.336
$ws.Cells.Item(123, 8).Value = 36486.668
$ws.Cells.Item(123, 10).Value = 36486.668
$ws.Cells.Item(123, 12).Value = 36486.668
$ws.Cells.Item(123, 14).Value = -46286.668
$ws.Cells.Item(125, 8).Value = 39435
$ws.Cells.Item(125, 10).Value = 39435
$ws.Cells.Item(125, 12).Value = 39435
$ws.Cells.Item(125, 14).Value = -49275
$ws.Cells.Item(132, 8).Value = 6064451
$ws.Cells.Item(132, 10).Value = 11907820
$ws.Cells.Item(132, 12).Value = 35723460
$ws.Cells.Item(132, 14).Value = -35728520
$ws.Cells.Item(134, 8).Value = 70607.336
$ws.Cells.Item(134, 10).Value = 70607.336
$ws.Cells.Item(134, 12).Value = 211822.008
$ws.Cells.Item(134, 14).Value = -216892.008
$ws.Cells.Item(136, 8).Value = 5770.086
$ws.Cells.Item(136, 9).Value = 6326.8335
$ws.Cells.Item(136, 10).Value = 5180.5884
$ws.Cells.Item(136, 11).Value = 18980.5005
$ws.Cells.Item(136, 12).Value = 15541.7652
$ws.Cells.Item(136, 13).Value = -16430.5005
$ws.Cells.Item(136, 14).Value = -20641.7652
$ws.Cells.Item(141, 8).Value = 47910
$ws.Cells.Item(141, 10).Value = 47910
$ws.Cells.Item(141, 12).Value = 47910
$ws.Cells.Item(141, 14).Value = -58270
